# CH03.docx first-pass transcription touch-up (c56-0020)
#
# 1) Narrow the first (marginal-note) column of the opening table from
#    2160 dxa to 1713 dxa (both the table grid and the cell width).
# 2) Re-join the split run "become a stu" + "dent" into a single run
#    "become a student" (the mid-word split was an artifact of an inline
#    _GoBack bookmark sitting between the two runs).
# 3) Put the (hidden/system) _GoBack bookmark back at the very start of
#    that first cell's content instead of mid-word.

$d = $word.ActiveDocument

# --- 1. Resize the first column of the first (running-head) table -------
$tbl = $d.Tables.Item(1)
$col1 = $tbl.Columns.Item(1)
$col1.Width = 85.65   # points == 1713 dxa (20 dxa per point)

# --- 2. Merge "become a stu" / "dent" back into one run ------------------
# Re-typing the already-correct phrase over itself forces Word to
# collapse the two runs (and drops the bookmark that used to live
# between them) into a single contiguous run of text.
$found = $d.Content.Find.Execute("become a student", $true, $false, $false,
                                  $false, $false, $true, 1, $false,
                                  "become a student", 2)

# --- 3. Restore the _GoBack bookmark at the start of the first cell ------
$cell1 = $tbl.Cell(1, 1)
$startPos = $cell1.Range.Start + 1
$bmRange = $d.Range($startPos, $startPos)
$bmAdded = $d.Bookmarks.Add("_GoBack", $bmRange)
